# Chile primera-division 2023: fix mis-ordered match rows and append
# five newly scraped matches at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Several row pairs had their match-detail columns (F..V) swapped
#    with one another while the leading index/date columns (A..E)
#    stayed put. Swap the F:V content back into the correct rows.
# ---------------------------------------------------------------------
function Swap-MatchColumns {
    param($ws, $rowA, $rowB)

    $rngA = $ws.Range($ws.Cells.Item($rowA, 6), $ws.Cells.Item($rowA, 22))
    $rngB = $ws.Range($ws.Cells.Item($rowB, 6), $ws.Cells.Item($rowB, 22))

    $valsA = $rngA.Value2
    $valsB = $rngB.Value2

    $rngA.Value2 = $valsB
    $rngB.Value2 = $valsA
}

Swap-MatchColumns $ws 26 27
Swap-MatchColumns $ws 45 46
Swap-MatchColumns $ws 63 64
Swap-MatchColumns $ws 67 68
Swap-MatchColumns $ws 96 97

# ---------------------------------------------------------------------
# 2) Append five new match rows (184-188) at the end of the sheet,
#    matching the formatting of the last existing row (183).
# ---------------------------------------------------------------------
$newRows = @(
    @(183, "chile", "primera-division", "2023", 45192.72916666666, "Palestino",    2, "U. Espanola", 1, 1.81, "15/09/2023 16:42", 1.95, "23/09/2023 17:29", 3.74, "15/09/2023 16:42", 3.73, "23/09/2023 17:29", 4.12, "15/09/2023 16:42", 3.91, "23/09/2023 17:29", "https://www.betexplorer.com/football/chile/primera-division/palestino-u-espanola/r1F59b2R/"),
    @(184, "chile", "primera-division", "2023", 45192.83333333334, "Curico Unido", 0, "Nublense",    3, 2.33, "15/09/2023 19:42", 3.47, "23/09/2023 19:58", 3.39, "15/09/2023 19:42", 3.27, "23/09/2023 19:58", 2.97, "15/09/2023 19:42", 2.28, "23/09/2023 19:58", "https://www.betexplorer.com/football/chile/primera-division/curico-unido-nublense/txet3due/"),
    @(185, "chile", "primera-division", "2023", 45192.9375,         "Colo Colo",    6, "Cobresal",    0, 1.66, "15/09/2023 21:42", 1.77, "23/09/2023 22:27", 4.17, "15/09/2023 21:42", 3.82, "23/09/2023 22:27", 4.87, "15/09/2023 21:42", 4.77, "23/09/2023 22:27", "https://www.betexplorer.com/football/chile/primera-division/colo-colo-cobresal/jgcbBKXE/"),
    @(186, "chile", "primera-division", "2023", 45193.04166666666, "Everton",      1, "Huachipato",  2, 2.06, "16/09/2023 00:43", 2.25, "24/09/2023 00:57", 3.55, "16/09/2023 00:43", 3.42, "24/09/2023 00:57", 3.41, "16/09/2023 00:43", 3.36, "24/09/2023 00:57", "https://www.betexplorer.com/football/chile/primera-division/everton-huachipato/zafp2Gf1/"),
    @(187, "chile", "primera-division", "2023", 45193.83333333334, "Copiapo",      3, "U. De Chile", 1, 2.71, "16/09/2023 19:43", 3.82, "24/09/2023 19:58", 3.32, "16/09/2023 19:43", 3.63, "24/09/2023 19:58", 2.69, "16/09/2023 19:43", 2.01, "24/09/2023 19:58", "https://www.betexplorer.com/football/chile/primera-division/copiapo-u-de-chile/8U7gC0I8/")
)

$lastRow = 183
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i

    # Copy A:V formatting from the previous last data row so the new
    # row gets the same styles (bold index column, date-formatted
    # match-date column, etc.) without disturbing the sheet dimension.
    $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 22)).Copy()
    $ws.Range($ws.Cells.Item($targetRow, 1), $ws.Cells.Item($targetRow, 22)).PasteSpecial(-4122)

    $row = $newRows[$i]
    $ws.Cells.Item($targetRow, 1).Value2  = $row[0]
    $ws.Cells.Item($targetRow, 2).Value2  = $row[1]
    $ws.Cells.Item($targetRow, 3).Value2  = $row[2]
    $ws.Cells.Item($targetRow, 5).Value2  = $row[4]
    $ws.Cells.Item($targetRow, 6).Value2  = $row[5]
    $ws.Cells.Item($targetRow, 7).Value2  = $row[6]
    $ws.Cells.Item($targetRow, 8).Value2  = $row[7]
    $ws.Cells.Item($targetRow, 9).Value2  = $row[8]
    $ws.Cells.Item($targetRow, 10).Value2 = $row[9]
    $ws.Cells.Item($targetRow, 11).Value2 = $row[10]
    $ws.Cells.Item($targetRow, 12).Value2 = $row[11]
    $ws.Cells.Item($targetRow, 13).Value2 = $row[12]
    $ws.Cells.Item($targetRow, 14).Value2 = $row[13]
    $ws.Cells.Item($targetRow, 15).Value2 = $row[14]
    $ws.Cells.Item($targetRow, 16).Value2 = $row[15]
    $ws.Cells.Item($targetRow, 17).Value2 = $row[16]
    $ws.Cells.Item($targetRow, 18).Value2 = $row[17]
    $ws.Cells.Item($targetRow, 19).Value2 = $row[18]
    $ws.Cells.Item($targetRow, 20).Value2 = $row[19]
    $ws.Cells.Item($targetRow, 21).Value2 = $row[20]
    $ws.Cells.Item($targetRow, 22).Value2 = $row[21]

    # Column D ("temporada") holds the purely-numeric-looking text
    # "2023". Assigning that via Value2 gets auto-coerced to a number
    # by the COM variant marshalling, so instead copy the *value* of
    # the existing, already-text D183 cell into this row's D cell -
    # a value-only paste preserves the original string type and adds
    # no style.
    $ws.Cells.Item($lastRow, 4).Copy()
    $ws.Cells.Item($targetRow, 4).PasteSpecial(-4163)  # xlPasteValues
}

Write-Output "done"
